$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44964
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2250
$ws.Range("P2").Value = 750

# Row 3
$ws.Range("D3").Value = 44883
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("P3").Value = 633

# Row 5
$ws.Range("D5").Value = 44827
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2250
$ws.Range("P5").Value = 750

# Row 6
$ws.Range("D6").Value = 44951
$ws.Range("J6").Value = 800

# Row 8
$ws.Range("D8").Value = 44953
$ws.Range("J8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44970
$ws.Range("J9").Value = 800
$ws.Range("K9").Value = 2000
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = 2250
$ws.Range("P9").Value = 750

# Row 10
$ws.Range("D10").Value = 44685
$ws.Range("J10").Value = 400

# Row 11
$ws.Range("D11").Value = 44881
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 1900
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1950
$ws.Range("P11").Value = 650

# Row 12
$ws.Range("D12").Value = 44910
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 1800
$ws.Range("M12").Value = 1900
$ws.Range("P12").Value = 633

# Row 13
$ws.Range("D13").Value = 44911
$ws.Range("J13").Value = 700
$ws.Range("K13").Value = 1800
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 1900
$ws.Range("P13").Value = 633

# Row 14
$ws.Range("D14").Value = 44848
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 1500
$ws.Range("M14").Value = 1750
$ws.Range("P14").Value = 583
